$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells (AD1:AF1), copying the header style (bold, border,
# centered) from the existing header cell AC1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record (Wins / Losses / Ties) for every data row (2-40).
$ws.Range("AD2:AD40").Value = 89
$ws.Range("AE2:AE40").Value = 73
$ws.Range("AF2:AF40").Value = 0

Write-Host "done"
